$d = $word.ActiveDocument

# --- Step 1: the "7 oktober" entry was originally split across two runs
# (an artifact of how it was authored). Merge them into a single run by
# doing a Find/Replace across the run boundary with the identical text -
# Word naturally consolidates the replaced span into one run.
$oldBoundary = "hdmi. Ook is onze versie"
$found = $d.Content.Find.Execute($oldBoundary, $true, $false, $false, $false, $false, $true, 1, $false, $oldBoundary, 2)
if (-not $found) {
    throw "Could not find the pi-top/pi-4b run boundary text to merge."
}

# --- Step 2: locate the (now single-run) "7 oktober" paragraph so we can
# append the new diary entries right after it (and after its trailing
# page-break run, which stays put at the end of that paragraph).
$targetText = "pi-4b binnen gekregen"
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($targetText)) {
        $anchorPara = $i
        break
    }
}
if (-not $anchorPara) {
    throw "Could not locate the pi-4b paragraph."
}

$newEntries = @(
    "8 oktober:",
    "thuis gewerkt. Requirements en gemaakte keuzes beschreven.",
    "9 oktober:",
    "Requirements besproken met Projectbegeleider en aangepast. Code voor het veranderen van de taal binnen het programma geschreven.",
    "12 oktober:",
    "Code geschreven voor het wisselen tussen een geavanceerde en normale modus. In de geavanceerde modus kunnen bijvoorbeeld bitwise operations gebruikt worden. Deze zullen waarschijnlijk te moeilijk zijn (en te weinig echt nut hebben) om gebruikt te worden, maar kunnen nu voor de volledigheid wel gebruikt worden."
)

$anchorRange = $d.Paragraphs.Item($anchorPara).Range
$insertIndex = $anchorPara
foreach ($entry in $newEntries) {
    $insertIndex = $insertIndex + 1
    $anchorRange.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertIndex)
    $newPara.Range.Text = $entry
    $anchorRange = $newPara.Range
}
